$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45345345
$ws.Range("B1").Value = 56756756

$ws.Range("B1").Select()
